{"js": "// Update each two-digit multiplication problem's text with its new value.\n// Every \"old\" cell value below occurs exactly once in the document body,\n// so a literal (non-wildcard) search-and-replace per pair is unambiguous.\nconst replacements = [\n  { oldText: \"76\u00d776=\", newText: \"24\u00d781=\" },\n  { oldText: \"74\u00d748=\", newText: \"75\u00d765=\" },\n  { oldText: \"79\u00d712=\", newText: \"50\u00d717=\" },\n  { oldText: \"36\u00d739=\", newText: \"52\u00d733=\" },\n  { oldText: \"21\u00d744=\", newText: \"49\u00d714=\" },\n  { oldText: \"41\u00d781=\", newText: \"90\u00d779=\" },\n  { oldText: \"39\u00d733=\", newText: \"26\u00d728=\" },\n  { oldText: \"20\u00d722=\", newText: \"87\u00d717=\" },\n  { oldText: \"33\u00d773=\", newText: \"49\u00d736=\" },\n  { oldText: \"46\u00d723=\", newText: \"90\u00d728=\" },\n  { oldText: \"91\u00d732=\", newText: \"76\u00d791=\" },\n  { oldText: \"51\u00d735=\", newText: \"49\u00d755=\" },\n  { oldText: \"53\u00d776=\", newText: \"75\u00d742=\" },\n  { oldText: \"25\u00d733=\", newText: \"57\u00d745=\" },\n  { oldText: \"39\u00d736=\", newText: \"49\u00d767=\" },\n  { oldText: \"61\u00d713=\", newText: \"81\u00d792=\" },\n  { oldText: \"27\u00d765=\", newText: \"53\u00d746=\" },\n  { oldText: \"67\u00d712=\", newText: \"96\u00d756=\" },\n  { oldText: \"18\u00d797=\", newText: \"53\u00d734=\" },\n  { oldText: \"87\u00d755=\", newText: \"80\u00d772=\" },\n  { oldText: \"94\u00d791=\", newText: \"53\u00d762=\" },\n  { oldText: \"34\u00d798=\", newText: \"21\u00d739=\" },\n  { oldText: \"27\u00d793=\", newText: \"44\u00d771=\" },\n  { oldText: \"17\u00d727=\", newText: \"85\u00d745=\" },\n  { oldText: \"83\u00d742=\", newText: \"80\u00d722=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication problem's text with its new value.\n# Every \"old\" cell value below occurs exactly once in the document, so\n# Find/Replace (wdReplaceAll = 2) per pair is unambiguous and idempotent.\n$pairs = @(\n    @{ Old = \"76\u00d776=\"; New = \"24\u00d781=\" },\n    @{ Old = \"74\u00d748=\"; New = \"75\u00d765=\" },\n    @{ Old = \"79\u00d712=\"; New = \"50\u00d717=\" },\n    @{ Old = \"36\u00d739=\"; New = \"52\u00d733=\" },\n    @{ Old = \"21\u00d744=\"; New = \"49\u00d714=\" },\n    @{ Old = \"41\u00d781=\"; New = \"90\u00d779=\" },\n    @{ Old = \"39\u00d733=\"; New = \"26\u00d728=\" },\n    @{ Old = \"20\u00d722=\"; New = \"87\u00d717=\" },\n    @{ Old = \"33\u00d773=\"; New = \"49\u00d736=\" },\n    @{ Old = \"46\u00d723=\"; New = \"90\u00d728=\" },\n    @{ Old = \"91\u00d732=\"; New = \"76\u00d791=\" },\n    @{ Old = \"51\u00d735=\"; New = \"49\u00d755=\" },\n    @{ Old = \"53\u00d776=\"; New = \"75\u00d742=\" },\n    @{ Old = \"25\u00d733=\"; New = \"57\u00d745=\" },\n    @{ Old = \"39\u00d736=\"; New = \"49\u00d767=\" },\n    @{ Old = \"61\u00d713=\"; New = \"81\u00d792=\" },\n    @{ Old = \"27\u00d765=\"; New = \"53\u00d746=\" },\n    @{ Old = \"67\u00d712=\"; New = \"96\u00d756=\" },\n    @{ Old = \"18\u00d797=\"; New = \"53\u00d734=\" },\n    @{ Old = \"87\u00d755=\"; New = \"80\u00d772=\" },\n    @{ Old = \"94\u00d791=\"; New = \"53\u00d762=\" },\n    @{ Old = \"34\u00d798=\"; New = \"21\u00d739=\" },\n    @{ Old = \"27\u00d793=\"; New = \"44\u00d771=\" },\n    @{ Old = \"17\u00d727=\"; New = \"85\u00d745=\" },\n    @{ Old = \"83\u00d742=\"; New = \"80\u00d722=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
